$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new DOI citation row (row 4)
$ws.Range("A4").Value = 10.14443
$ws.Range("B4").Value = 9999

# Update the active cell selection to A5 (as in the saved file)
$ws.Range("A5").Select()
